# Append three new case records (rows 61-63) to the bottom of the case
# data table on Sheet1, matching the source export format:
#   A..G  = text fields (case #, name, charge, statute, degree, plea, finding)
#   H     = numeric (amount), 0 here
#   I     = text "0" (kept as text even though it looks numeric)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @("21CRB01268", "Bunner", "POSSESSION DRUG PARAPHERNALIA", "2925.14(C)", "M4", "Guilty", "Guilty", 0, "0"),
    @("21CRB01291", "Bunner", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM", "No Contest", "Guilty", 0, "0"),
    @("21CRB01291", "Bunner", "No Operator License - Never Held", "4510.12(C)(1)", "Unclassified Misdemeanor", "No Contest", "Guilty", 0, "0")
)

$startRow = 61
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    # Column I holds "0" as TEXT, not a number. A leading apostrophe forces
    # Excel to store it as text (quote-prefixed) instead of coercing it to
    # a numeric value; resetting the style back to "Normal" afterwards
    # drops the quote-prefix formatting flag while the cell keeps its text
    # type, so no stray number format is left behind on the cell.
    $ws.Cells.Item($r, 9).Value = "'" + $row[8]
    $ws.Cells.Item($r, 9).Style = "Normal"
}
